# Apply "updated 4.0 files and mdl" edits to the Maximum Capacity Factor workbook
$wb = $excel.ActiveWorkbook

# --- About sheet: update the date stamp in C1 (serial 45320 -> 45392) ---
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = [DateTime]::FromOADate(45392)

# --- MCF sheet: bump capacity factors from 0.85/0.95 up to 1 ---
$mcf = $wb.Worksheets.Item("MCF")

$cellsToOne = @("B2","B3","B4","B6","B10","B11","B12","B13","B14","B16","B17","B18")
foreach ($addr in $cellsToOne) {
    $mcf.Range($addr).Value = 1
}

# Cells with formulas referencing the cells above will recalculate automatically
# (B19=B2, B20=B4, B21=B10, B22=B14, B24=B4, B25=B4)

$excel.Calculate()

# --- Restore the active sheet selection to match the saved workbook state ---
$mcf.Activate()
$mcf.Range("B17").Select()
